$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: model_type -> type
$ws.Range("H1").Value = "type"

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count

for ($r = 2; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 8)  # column H is 8th column
    $val = $cell.Value2
    if ($val -eq "link" -or $val -eq "source" -or $val -eq "sink") {
        $cell.Value = "road"
    }
}
